$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.490922689437866
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 4.733295440673828
$ws.Range("D1").Value = 2.222400665283203
$ws.Range("E1").Value = 1.734768033027649
